$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "4510 Insulation;"; New = "Design: 4510 Insulation;" },
    @{ Old = "Yes"; New = "Design: Yes" },
    @{ Old = "It's a bit unclear if A-15 fire insulation needs overlaps. In 518 the shipowner wanted to have overlaps around pipes (penetrating A-15 bulkhead)."; New = "Design: It's a bit unclear if A-15 fire insulation needs overlaps. In 518 the shipowner wanted to have overlaps around pipes (penetrating A-15 bulkhead)." },
    @{ Old = "Budget was exceeded, material prices had risen a lot from 516 (from year 2021)"; New = "Design: Budget was exceeded, material prices had risen a lot from 516 (from year 2021)" },
    @{ Old = "When insulation contractor suggest's changes in drawing (scheme of insulation), it's important to check the cost effect for shipyard."; New = "Design: When insulation contractor suggest's changes in drawing (scheme of insulation), it's important to check the cost effect for shipyard." },
    @{ Old = "Co-operation was good between contractor, production insulation foreman and design."; New = "Design: Co-operation was good between contractor, production insulation foreman and design." },
    @{ Old = "No difference"; New = "Design: No difference" },
    @{ Old = "Pins for insulation were installed too late - was planned to be installed in block factory, but quite much was installed in blocks or onboard after painting. "; New = "Design: Pins for insulation were installed too late - was planned to be installed in block factory, but quite much was installed in blocks or onboard after painting. " }
)

foreach ($r in $replacements) {
    $rng = $d.Content
    # Locate the exact text (whole-word match, case-sensitive) without
    # letting Find perform the replacement itself, since Find's built-in
    # replace pass runs text through AutoCorrect (straight quotes become
    # curly "smart" quotes). wdFindStop / wdReplaceNone = no substitution.
    $found = $rng.Find.Execute($r.Old, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Text = $r.New
    }
}
